# "modif algorithme creation des nouvelles tables"
#
# The sheet flags, for every observation row, which "new tables"
# (SGI / SGI+DICR / SGResp / GD+PARAM / HAEM / RESIST -> columns D..I)
# apply. Until now every flagged column simply held a "1". The new
# table-creation algorithm instead stores a fixed per-column weight
# (D=1, E=2, F=3, G=4, H=5, I=6) so downstream code can tell which
# table produced the match. Cells that were empty stay empty; only
# cells that already carried a flag get rewritten with their column's
# weight.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index (1-based) -> new weight value.
$colWeights = @{
    4 = 1   # D: SGI
    5 = 2   # E: SGI + DICR
    6 = 3   # F: SGResp
    7 = 4   # G: GD + PARAM
    8 = 5   # H: HAEM
    9 = 6   # I: RESIST
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 38 }

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $colWeights.Keys) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null -and $v -ne "") {
            $cell.Value = $colWeights[$c]
        }
    }
}

# Reflect the view/selection state captured in the saved workbook: the
# window had scrolled so row 29 is the first visible row, and the
# active selection is I2:I39 (anchored at I2).
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("I2:I39").Select()
